$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.438.74'
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.744.30'
$ws.Range("E3").Value = '  -3.45%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.95'
$ws.Range("E5").Value = '  -4.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4239'
$ws.Range("E7").Value = '  -8.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3585'
$ws.Range("E8").Value = '  -3.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.48'
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07428'
$ws.Range("E10").Value = '  -3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  -3.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.112'
$ws.Range("E14").Value = '  -3.68%  '
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.743.68'
$ws.Range("E16").Value = '  -3.45%  '
$ws.Range("E17").Value = '  -3.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.16'
$ws.Range("E18").Value = '  +7.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06116'
$ws.Range("E19").Value = '  -8.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.88'
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.102'
$ws.Range("E22").Value = '  -4.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5261'
$ws.Range("E23").Value = '  -6.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.474.24'
$ws.Range("E24").Value = '  -2.76%  '
$ws.Range("E25").Value = '  -3.46%  '
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.96'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.377'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.940.92'
$ws.Range("E30").Value = '  -3.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.90'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.196'
$ws.Range("E32").Value = '  -4.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.672'
$ws.Range("E33").Value = '  -3.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09117'
$ws.Range("E34").Value = '  -4.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.625'
$ws.Range("E35").Value = '  -10.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.67'
$ws.Range("E36").Value = '  +4.42%  '
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2137'
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.083'
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06055'
$ws.Range("E40").Value = '  -4.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6387'
$ws.Range("E41").Value = '  -3.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.189'
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.430'
$ws.Range("E43").Value = '  -4.58%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.892'
$ws.Range("E45").Value = '  -4.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.69'
$ws.Range("E46").Value = '  -4.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.714'
$ws.Range("E47").Value = '  -2.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5867'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.07'
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("E50").Value = '  -5.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06840'
$ws.Range("E51").Value = '  -4.51%  '
